$wb = $excel.ActiveWorkbook

# --- Sheet2: append new operator rows (8-18) ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$data = @(
    @("BETWEEN", "Between minimum and maximum values"),
    @("NOT BETWEEN", "Outside minimum and maximum values"),
    @("NOT", "It is not the case that …"),
    @("IN", "Matches one of a list of values"),
    @("NOT IN", "Does not match any of a list of values"),
    @("IS NULL", "Has no value"),
    @("IS NOT NULL", "Has some value"),
    @("LIKE", "Matches pattern"),
    @("NOT LIKE", "Does not match pattern"),
    @("ANY/SOME", "Condition applies to any of a list of values"),
    @("ALL", "Condition applies to every value in a list of values ")
)

$row = 8
foreach ($pair in $data) {
    $ws2.Cells.Item($row, 1).Value = $pair[0]
    $ws2.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# --- Update active cell selections on each sheet ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("H17").Select()

$ws2.Range("E14").Select()

# Sheet1 was (and should remain) the active/visible tab; re-activate it last so the
# workbook re-opens on Sheet1 while each sheet keeps its own selection above.
$ws1.Activate()
